# feat: add 2022-Q4 data
#
# 1) Insert a new worksheet "2022-Q4" right after "总计" (i.e. immediately
#    before the existing "2022-Q3" sheet), populated with the fund-holding
#    detail rows for the new quarter.
# 2) Insert a new leading row in the "总计" summary sheet for 2022-Q4,
#    pushing the existing quarters down by one row.

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# Helper style references taken from the existing "总计" sheet: a header /
# index-column cell (bold, thin border, centered) and a plain data cell
# (no special formatting). We clone these via Copy + PasteSpecial(formats)
# so the new cells land on the exact same style entries Excel already uses
# in this workbook, instead of synthesizing new ones.
# ---------------------------------------------------------------------------
$headerStyleSrc = $summary.Range("B1")
$plainStyleSrc = $summary.Range("B3")

# ---------------------------------------------------------------------------
# 1) Insert the new sheet before the current second sheet ("2022-Q3").
# ---------------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# 2) Populate the new "2022-Q4" sheet header row (B1:H1).
# ---------------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 2).Value = $headers[$i]
}
$headerStyleSrc.Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Populate the fund-holding detail rows.
#    Column A: plain 0-based running index (number).
#    Columns B-G: kept as TEXT (matches the scraped source, which stores
#      fund codes / sizes / percentages as inline strings, not numbers).
#    Column H: numeric position rank.
# ---------------------------------------------------------------------------
$rows = @(
    @("517160", "南方中证长江保护主题ETF", "16.77", "99.26", "2.50", "0.4192", 4),
    @("517330", "易方达中证长江保护主题ETF", "16.47", "99.39", "2.51", "0.4134", 4),
    @("159610", "景顺长城中证500增强策略ETF", "6.09", "98.72", "1.41", "0.0859", 3),
    @("501030", "汇添富中证环境治理指数（LOF）A", "2.99", "92.14", "1.94", "0.0580", 6),
    @("164908", "交银施罗德中证环境治理指数（LOF）", "1.55", "93.92", "1.99", "0.0308", 6),
    @("501031", "汇添富中证环境治理指数（LOF）C", "1.38", "92.14", "1.94", "0.0268", 6),
    @("012879", "中信建投量化精选6个月持有期混合C", "3.12", "70.35", "0.74", "0.0231", 9),
    @("012878", "中信建投量化精选6个月持有期混合A", "1.59", "70.35", "0.74", "0.0118", 9),
    @("005260", "银华稳健增利灵活配置混合A", "0.33", "91.18", "0.66", "0.0022", 8),
    @("013413", "交银施罗德中证环境治理指数（LOF）C", "0.11", "93.92", "1.99", "0.0022", 6),
    @("005261", "银华稳健增利灵活配置混合C", "0.21", "91.18", "0.66", "0.0014", 8)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $excelRow = $r + 2

    $newSheet.Cells.Item($excelRow, 1).Value = $r

    for ($c = 0; $c -lt 6; $c++) {
        $cell = $newSheet.Cells.Item($excelRow, $c + 2)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$c]
    }

    $newSheet.Cells.Item($excelRow, 8).Value = $row[6]
}

$headerStyleSrc.Copy()
$newSheet.Range("A2:A12").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) Update the "总计" (summary) sheet: insert a new leading data row for
#    2022-Q4, shifting the existing quarters down by one row. Column A is a
#    plain 0-based running index over the data rows, so it must be
#    renumbered after the shift (not merely carried along).
# ---------------------------------------------------------------------------
$summary.Rows(2).Insert()

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 11
$summary.Cells.Item(2, 4).Value = 1.07

$headerStyleSrc.Copy()
$summary.Range("A2").PasteSpecial(-4122)
$plainStyleSrc.Copy()
$summary.Range("B2:D2").PasteSpecial(-4122)

for ($r = 3; $r -le 7; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}
